# Applies the "Updated cryptos list" data refresh to Sheet1.
# Column D values are forced to Text via NumberFormat "@" before the
# assignment (then restored) so numeric-looking strings such as
# "0.999" or "54.00" are not silently converted to real numbers by Excel.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $origFormat = $range.NumberFormat
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.NumberFormat = $origFormat
}

Set-TextValue $ws.Range("D2") "61.930.06"
$ws.Range("E2").Value = "  -0.35%  "

Set-TextValue $ws.Range("D3") "3.406.69"
$ws.Range("E3").Value = "  -0.65%  "

Set-TextValue $ws.Range("D4") "0.999"
$ws.Range("E4").Value = "  -0.13%  "

Set-TextValue $ws.Range("D5") "408.14"
$ws.Range("E5").Value = "  +0.03%  "

Set-TextValue $ws.Range("D6") "128.79"
$ws.Range("E6").Value = "  -3.47%  "

Set-TextValue $ws.Range("D7") "0.635"
$ws.Range("E7").Value = "  +7.07%  "

Set-TextValue $ws.Range("D9") "0.729"
$ws.Range("E9").Value = "  +7.41%  "

Set-TextValue $ws.Range("D10") "0.143"
$ws.Range("E10").Value = "  +16.36%  "

Set-TextValue $ws.Range("D11") "42.35"
$ws.Range("E11").Value = "  -0.23%  "

$ws.Range("B12").Value = "ShibaInu"
$ws.Range("C12").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue $ws.Range("D12") "0.0000216"
$ws.Range("E12").Value = "  +64.68%  "

$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue $ws.Range("D13") "0.140"
$ws.Range("E13").Value = "  -0.57%  "

$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextValue $ws.Range("D14") "3.949.13"
$ws.Range("E14").Value = "  -0.37%  "

Set-TextValue $ws.Range("D15") "8.89"
$ws.Range("E15").Value = "  +5.43%  "

Set-TextValue $ws.Range("D16") "20.83"
$ws.Range("E16").Value = "  +4.60%  "

Set-TextValue $ws.Range("D17") "3.400.28"
$ws.Range("E17").Value = "  -0.37%  "

Set-TextValue $ws.Range("D18") "12.05"
$ws.Range("E18").Value = "  +9.14%  "

Set-TextValue $ws.Range("D19") "1.07"
$ws.Range("E19").Value = "  +4.65%  "

Set-TextValue $ws.Range("D20") "61.799.06"
$ws.Range("E20").Value = "  -0.53%  "

Set-TextValue $ws.Range("D21") "404.64"
$ws.Range("E21").Value = "  +28.29%  "

Set-TextValue $ws.Range("D22") "89.21"
$ws.Range("E22").Value = "  +5.08%  "

Set-TextValue $ws.Range("D23") "3.18"
$ws.Range("E23").Value = "  -1.09%  "

Set-TextValue $ws.Range("D24") "13.06"
$ws.Range("E24").Value = "  +1.91%  "

Set-TextValue $ws.Range("D25") "3.22"
$ws.Range("E25").Value = "  +3.13%  "

Set-TextValue $ws.Range("D26") "32.75"
$ws.Range("E26").Value = "  +10.15%  "

$ws.Range("B27").Value = "Filecoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws.Range("D27") "8.60"
$ws.Range("E27").Value = "  +4.37%  "

$ws.Range("B28").Value = "LEO"
$ws.Range("C28").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue $ws.Range("D28") "4.79"
$ws.Range("E28").Value = "  -0.07%  "

Set-TextValue $ws.Range("D29") "7.60"
$ws.Range("E29").Value = "  -1.66%  "

$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue $ws.Range("D30") "2.74"
$ws.Range("E30").Value = "  -0.50%  "

$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws.Range("D31") "0.117"
$ws.Range("E31").Value = "  +1.19%  "

$ws.Range("E32").Value = "  -1.92%  "

Set-TextValue $ws.Range("D33") "11.84"
$ws.Range("E33").Value = "  +3.89%  "

Set-TextValue $ws.Range("D34") "43.06"
$ws.Range("E34").Value = "  +0.33%  "

Set-TextValue $ws.Range("D36") "0.0496"
$ws.Range("E36").Value = "  +2.16%  "

Set-TextValue $ws.Range("D37") "54.00"
$ws.Range("E37").Value = "  +3.44%  "

Set-TextValue $ws.Range("D38") "0.997"
$ws.Range("E38").Value = "  -0.13%  "

$ws.Range("E39").Value = "  -2.88%  "

Set-TextValue $ws.Range("D40") "0.133"
$ws.Range("E40").Value = "  +6.35%  "

Set-TextValue $ws.Range("D41") "2.92"
$ws.Range("E41").Value = "  -2.53%  "

Set-TextValue $ws.Range("D42") "0.311"
$ws.Range("E42").Value = "  +5.84%  "

Set-TextValue $ws.Range("D43") "140.78"
$ws.Range("E43").Value = "  +2.17%  "

Set-TextValue $ws.Range("D44") "1.97"
$ws.Range("E44").Value = "  -1.93%  "

Set-TextValue $ws.Range("D45") "4.06"
$ws.Range("E45").Value = "  +1.63%  "

$ws.Range("E46").Value = "  +8.51%  "

Set-TextValue $ws.Range("D47") "16.63"
$ws.Range("E47").Value = "  -1.02%  "

Set-TextValue $ws.Range("D48") "21.73"
$ws.Range("E48").Value = "  +1.27%  "

Set-TextValue $ws.Range("D49") "2.113.72"
$ws.Range("E49").Value = "  -0.83%  "

$ws.Range("E50").Value = "  +4.17%  "

Set-TextValue $ws.Range("D51") "0.131"
$ws.Range("E51").Value = "  +15.28%  "
